$wb = $excel.ActiveWorkbook
$wsDemand = $wb.Worksheets.Item("Demand_vs_Served")
$wsGap = $wb.Worksheets.Item("Gap")
$wsCap = $wb.Worksheets.Item("Capacity_Utilization")

# --- Demand_vs_Served sheet: row 97 (last slot) had capacity drop to 0 ---
$wsDemand.Range("D97").Value = 0
$wsDemand.Range("E97").Value = 10

# --- Gap sheet: row 97 gap_tasks updated to reflect unmet demand ---
$wsGap.Range("C97").Value = 10

# --- Capacity_Utilization sheet: capacity_tasks and utilization recalculated per slot ---
$capChanges = @{
    "C2" = 6
    "E2" = 0.3333333333333333
    "C4" = 6
    "E4" = 1
    "C5" = 12
    "E5" = 0.8333333333333334
    "C6" = 12
    "E6" = 0.75
    "C7" = 18
    "E7" = 0.5
    "C8" = 12
    "E8" = 0.6666666666666666
    "C9" = 6
    "E9" = 0.5
    "C10" = 6
    "E10" = 0.6666666666666666
    "C11" = 6
    "E11" = 1
    "C12" = 6
    "E12" = 0.5
    "C13" = 12
    "E13" = 0.6666666666666666
    "C15" = 12
    "E15" = 0.5833333333333334
    "C16" = 6
    "E16" = 0.3333333333333333
    "C17" = 12
    "E17" = 0.6666666666666666
    "C18" = 6
    "E18" = 0.3333333333333333
    "C19" = 6
    "E19" = 0.8333333333333334
    "C20" = 12
    "E20" = 0.5833333333333334
    "C21" = 12
    "E21" = 0.75
    "C22" = 12
    "E22" = 0.75
    "C24" = 6
    "E24" = 0.6666666666666666
    "C25" = 6
    "E25" = 1
    "C26" = 12
    "E26" = 0.6666666666666666
    "C28" = 12
    "E28" = 0.5833333333333334
    "C29" = 12
    "E29" = 0.8333333333333334
    "C31" = 6
    "E31" = 0.5
    "C32" = 6
    "E32" = 1
    "C33" = 12
    "E33" = 0.75
    "C34" = 12
    "E34" = 0.5833333333333334
    "C35" = 12
    "E35" = 0.75
    "C36" = 12
    "E36" = 0.75
    "C37" = 12
    "E37" = 0.6666666666666666
    "C38" = 12
    "E38" = 0.75
    "C39" = 12
    "E39" = 0.6666666666666666
    "C40" = 12
    "E40" = 0.75
    "C42" = 12
    "E42" = 0.8333333333333334
    "C43" = 6
    "E43" = 0.3333333333333333
    "C44" = 6
    "E44" = 1
    "C45" = 12
    "E45" = 0.75
    "C46" = 12
    "E46" = 0.6666666666666666
    "C47" = 6
    "E47" = 0.5
    "C48" = 6
    "E48" = 1
    "C49" = 6
    "E49" = 0.8333333333333334
    "C50" = 12
    "E50" = 0.75
    "C52" = 6
    "E52" = 0.8333333333333334
    "C54" = 6
    "E54" = 0.8333333333333334
    "C55" = 12
    "E55" = 0.6666666666666666
    "C56" = 6
    "E56" = 0.8333333333333334
    "C57" = 6
    "E57" = 0.3333333333333333
    "C58" = 6
    "E58" = 0.3333333333333333
    "C60" = 6
    "E60" = 0.5
    "C61" = 6
    "E61" = 1
    "C62" = 6
    "E62" = 1
    "C63" = 6
    "E63" = 0.6666666666666666
    "C64" = 6
    "E64" = 0.6666666666666666
    "C65" = 12
    "E65" = 0.5833333333333334
    "C66" = 6
    "E66" = 1
    "C67" = 12
    "E67" = 0.75
    "C68" = 12
    "E68" = 0.5833333333333334
    "C69" = 6
    "E69" = 1
    "C70" = 6
    "E70" = 0.6666666666666666
    "C71" = 6
    "E71" = 0.5
    "C72" = 6
    "E72" = 0.5
    "C73" = 6
    "E73" = 0.3333333333333333
    "C74" = 6
    "E74" = 0.6666666666666666
    "C75" = 12
    "E75" = 0.5833333333333334
    "C76" = 18
    "E76" = 0.5
    "C78" = 6
    "E78" = 1
    "C80" = 6
    "E80" = 1
    "C81" = 12
    "E81" = 0.75
    "C82" = 12
    "E82" = 0.75
    "C83" = 6
    "E83" = 1
    "C84" = 12
    "E84" = 0.75
    "C86" = 12
    "E86" = 0.8333333333333334
    "C88" = 6
    "E88" = 1
    "C89" = 6
    "E89" = 1
    "C90" = 6
    "E90" = 0.5
    "C91" = 6
    "E91" = 1
    "C92" = 6
    "E92" = 0.3333333333333333
    "C93" = 6
    "E93" = 0.6666666666666666
    "C94" = 12
    "E94" = 0.75
    "C95" = 12
    "E95" = 0.6666666666666666
    "C96" = 6
    "E96" = 0.6666666666666666
    "C97" = 0
    "D97" = 0
    "E97" = 0
}
foreach ($cell in $capChanges.Keys) {
    $wsCap.Range($cell).Value = $capChanges[$cell]
}

Write-Host "Applied timeseries_results.xlsx updates"
